# Commit: "add new code for multiple test data like one data rone then second data row"
#
# Data edits:
#  - MAIN_CONTROLLER!C2:            "Remote"       -> "local"
#  - DATASHEET!D3:                  "FOS.xlsx"     -> "FOS3.xlsx"
#  - MOBILE_CONFIGURATION!F3:       "One Plus 11R" -> "d4a4d1d2"
#  - MOBILE_CONFIGURATION!H3:       13             -> 11
#
# Plus the cursor/selection was left on MOBILE_CONFIGURATION!F7 (that sheet
# becomes the active tab), after also having visited MAIN_CONTROLLER!D2 and
# DATASHEET!D3 while editing.

$wb = $excel.ActiveWorkbook

# --- MAIN_CONTROLLER -------------------------------------------------------
$wsMain = $wb.Worksheets.Item("MAIN_CONTROLLER")
$wsMain.Activate()
$wsMain.Range("C2").Value = "local"
$wsMain.Range("D2").Select()

# --- DATASHEET ---------------------------------------------------------
$wsData = $wb.Worksheets.Item("DATASHEET")
$wsData.Activate()
$wsData.Range("D3").Value = "FOS3.xlsx"
$wsData.Range("D3").Select()

# --- MOBILE_CONFIGURATION ---------------------------------------------------
$wsMobile = $wb.Worksheets.Item("MOBILE_CONFIGURATION")
$wsMobile.Activate()
$wsMobile.Range("F3").Value = "d4a4d1d2"
$wsMobile.Range("H3").Value = 11
$wsMobile.Range("F7").Select()
